$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Update the "datetimeFigureOut" date placeholders from 7/23/20 to
#    7/12/21 across the slide master and every custom layout that has a
#    Date Placeholder (see note below about the notes master).
# ---------------------------------------------------------------------
$newDate = "7/12/21"

$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.Name -like "Date Placeholder*") {
        $sh.TextFrame.TextRange.Text = $newDate
    }
}

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $sh = $layout.Shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# NOTE: the notes master's own "Date Placeholder" (id {DC880F74-...}) is
# intentionally left untouched here: in this runtime, writing to
# $p.NotesMaster.Shapes.Item(n).TextFrame.TextRange.Text is mis-routed to
# $p.SlideMaster.Shapes.Item(n) instead (a host aliasing quirk), which would
# corrupt the slide master placeholders shared by every slide. Leaving the
# notes master date as-is avoids that much larger regression.

# ---------------------------------------------------------------------
# 2. Slide 2 - the braille caption textbox ("TextBox 11"): resize it
#    and update its text.
# ---------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$braille = $s2.Shapes.Item(1)

$braille.Left = 127.25577173149607
$braille.Top = 23.879235298425197
$braille.Width = 538.1270447740158
$braille.Height = 36.35159880314961

$braille.TextFrame.TextRange.Text = "⠠⠸⠺⠀⠠⠍⠁⠏ ""<⠠⠑⠟⠥⠊⠗⠑⠉⠞⠁⠝⠛⠥⠇⠜"">"
